# This document contains a single table with practice problems of the form
# "AA÷B=" laid out in a 5-column grid. We replace each problems text with a
# new one. Find/Replace in this runtime operates over the whole document
# regardless of the Range it is invoked on, so the replacements below are
# ordered to avoid a value that is freshly written by one replacement being
# re-matched by a later replacement that is still searching for that same text
# as its "old" value (this occurs once, for "20÷9=").
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Find.Execute("30÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "42÷7=", 2) | Out-Null
$t.Cell(1,2).Range.Find.Execute("88÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "27÷9=", 2) | Out-Null
$t.Cell(1,3).Range.Find.Execute("72÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "91÷2=", 2) | Out-Null
$t.Cell(1,4).Range.Find.Execute("32÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "98÷6=", 2) | Out-Null
$t.Cell(1,5).Range.Find.Execute("64÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "29÷7=", 2) | Out-Null
$t.Cell(5,1).Range.Find.Execute("12÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "90÷8=", 2) | Out-Null
$t.Cell(5,2).Range.Find.Execute("54÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "23÷5=", 2) | Out-Null
$t.Cell(5,3).Range.Find.Execute("32÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "89÷2=", 2) | Out-Null
$t.Cell(5,4).Range.Find.Execute("17÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "25÷9=", 2) | Out-Null
$t.Cell(5,5).Range.Find.Execute("19÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "55÷7=", 2) | Out-Null
$t.Cell(9,1).Range.Find.Execute("31÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "41÷7=", 2) | Out-Null
$t.Cell(9,2).Range.Find.Execute("14÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "77÷8=", 2) | Out-Null
$t.Cell(9,3).Range.Find.Execute("77÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "40÷2=", 2) | Out-Null
$t.Cell(9,4).Range.Find.Execute("91÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "50÷8=", 2) | Out-Null
$t.Cell(9,5).Range.Find.Execute("66÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "41÷4=", 2) | Out-Null
$t.Cell(13,2).Range.Find.Execute("24÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "32÷5=", 2) | Out-Null
$t.Cell(13,3).Range.Find.Execute("78÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "56÷8=", 2) | Out-Null
$t.Cell(13,4).Range.Find.Execute("96÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "34÷2=", 2) | Out-Null
$t.Cell(13,5).Range.Find.Execute("72÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷3=", 2) | Out-Null
$t.Cell(17,1).Range.Find.Execute("84÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "33÷6=", 2) | Out-Null
$t.Cell(17,2).Range.Find.Execute("12÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "77÷6=", 2) | Out-Null
$t.Cell(17,3).Range.Find.Execute("20÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "66÷2=", 2) | Out-Null
$t.Cell(17,4).Range.Find.Execute("94÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "28÷8=", 2) | Out-Null
$t.Cell(17,5).Range.Find.Execute("19÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 2) | Out-Null
$t.Cell(13,1).Range.Find.Execute("47÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "20÷9=", 2) | Out-Null
